# Update the "想去人数" (want-to-go count, column F) figures across the
# three sheets that carry per-event rows: 展览, 演出, 全部类型.
# (本地生活 has no numeric change in this commit.)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 42
$ws.Range("F3").Value  = 21311
$ws.Range("F4").Value  = 816
$ws.Range("F6").Value  = 1131
$ws.Range("F7").Value  = 29
$ws.Range("F8").Value  = 7945
$ws.Range("F9").Value  = 554
$ws.Range("F10").Value = 43
$ws.Range("F12").Value = 315
$ws.Range("F13").Value = 65
$ws.Range("F15").Value = 170
$ws.Range("F16").Value = 32
$ws.Range("F19").Value = 1360
$ws.Range("F20").Value = 537
$ws.Range("F22").Value = 708
$ws.Range("F24").Value = 81
$ws.Range("F25").Value = 85
$ws.Range("F26").Value = 353
$ws.Range("F27").Value = 1190
$ws.Range("F28").Value = 54
$ws.Range("F32").Value = 606
$ws.Range("F35").Value = 5063
$ws.Range("F36").Value = 35
$ws.Range("F38").Value = 46
$ws.Range("F40").Value = 13133
$ws.Range("F41").Value = 1368
$ws.Range("F42").Value = 137
$ws.Range("F43").Value = 53
$ws.Range("F45").Value = 309
$ws.Range("F46").Value = 437
$ws.Range("F47").Value = 4062
$ws.Range("F49").Value = 103

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 329

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 42
$ws.Range("F3").Value  = 21311
$ws.Range("F4").Value  = 816
$ws.Range("F5").Value  = 1131
$ws.Range("F6").Value  = 29
$ws.Range("F7").Value  = 7945
$ws.Range("F8").Value  = 554
$ws.Range("F9").Value  = 43
$ws.Range("F11").Value = 315
$ws.Range("F12").Value = 65
$ws.Range("F14").Value = 170
$ws.Range("F15").Value = 32
$ws.Range("F17").Value = 1360
$ws.Range("F18").Value = 537
$ws.Range("F20").Value = 708
$ws.Range("F22").Value = 81
$ws.Range("F23").Value = 85
$ws.Range("F24").Value = 353
$ws.Range("F25").Value = 1190
$ws.Range("F26").Value = 54
$ws.Range("F29").Value = 329
$ws.Range("F30").Value = 606
$ws.Range("F35").Value = 5063
$ws.Range("F36").Value = 35
$ws.Range("F38").Value = 46
$ws.Range("F40").Value = 13133
$ws.Range("F41").Value = 1368
$ws.Range("F42").Value = 137
$ws.Range("F43").Value = 53
$ws.Range("F45").Value = 309
$ws.Range("F46").Value = 437
$ws.Range("F47").Value = 4062
$ws.Range("F49").Value = 103
